$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1. Remove the _GoBack bookmark that currently sits right after the
#    "Zachary C. Raslan" title text.
# ------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# ------------------------------------------------------------------
# 2. Insert three new empty Heading2-styled paragraphs right before the
#    "Work Experience" heading (i.e. right after the empty paragraph
#    that follows the education table). The first of these new
#    paragraphs gets the _GoBack bookmark re-created in its new spot.
# ------------------------------------------------------------------

# Locate the "Work Experience" heading using Find.
$findRange = $d.Content
$findRange.Find.ClearFormatting()
$found = $findRange.Find.Execute("Work Experience", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) {
    throw "Could not locate 'Work Experience' heading"
}
$targetStart = $findRange.Start

# Map the found character position back to its paragraph index in the
# document's Paragraphs collection (more reliable across this runtime
# than deriving it from a sub-range's own Paragraphs collection).
$workExpIndex = -1
$count = $d.Paragraphs.Count
for ($i = 1; $i -le $count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Start -le $targetStart -and $p.Range.End -gt $targetStart) {
        $workExpIndex = $i
        break
    }
}
if ($workExpIndex -eq -1) {
    throw "Could not map 'Work Experience' text to a paragraph"
}

$targetPara = $d.Paragraphs.Item($workExpIndex - 1)

# Create an insertion point right after that paragraph's mark.
$targetPara.Range.InsertParagraphAfter()
$insertRange = $d.Paragraphs.Item($workExpIndex).Range

$pWithBookmark = '<w:p><w:pPr><w:pStyle w:val="Heading2"/><w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:hAnsiTheme="minorHAnsi"/><w:sz w:val="32"/></w:rPr></w:pPr><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p>'
$pPlain = '<w:p><w:pPr><w:pStyle w:val="Heading2"/><w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:hAnsiTheme="minorHAnsi"/><w:sz w:val="32"/></w:rPr></w:pPr></w:p>'

$bodyInner = $pWithBookmark + $pPlain + $pPlain
$packageXml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' + `
    '<w:body>' + $bodyInner + '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$null = $insertRange.InsertXML($packageXml)
